# Update Leve market-profit figures (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) across the Garuda_Profits sheets, per the scheduled market-price
# refresh run.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1170.6897
$ws.Range("I18").Value = 1090.7407
$ws.Range("J18").Value = 2250
$ws.Range("K18").Value = 1090.7407
$ws.Range("L18").Value = 2250
$ws.Range("M18").Value = -806.7407000000001
$ws.Range("N18").Value = -2818
# Row 28
$ws.Range("H28").Value = 1852.24
$ws.Range("I28").Value = 289.1875
$ws.Range("J28").Value = 4631
$ws.Range("K28").Value = 289.1875
$ws.Range("L28").Value = 4631
$ws.Range("M28").Value = 195.8125
$ws.Range("N28").Value = -5601
# Row 41
$ws.Range("H41").Value = 4988.8335
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 7433.25
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 7433.25
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -8313.25
# Row 62
$ws.Range("H62").Value = 1063.0769
$ws.Range("I62").Value = 1112.3334
$ws.Range("K62").Value = 1112.3334
$ws.Range("M62").Value = -488.3334
# Row 65
$ws.Range("H65").Value = 1063.0769
$ws.Range("I65").Value = 1112.3334
$ws.Range("K65").Value = 5561.666999999999
$ws.Range("M65").Value = -2441.666999999999
# Row 86
$ws.Range("H86").Value = 58883572
$ws.Range("I86").Value = 112372.664
$ws.Range("J86").Value = 125001170
$ws.Range("K86").Value = 112372.664
$ws.Range("L86").Value = 125001170
$ws.Range("M86").Value = -111249.664
$ws.Range("N86").Value = -125003416
# Row 89
$ws.Range("H89").Value = 58883572
$ws.Range("I89").Value = 112372.664
$ws.Range("J89").Value = 125001170
$ws.Range("K89").Value = 561863.3200000001
$ws.Range("L89").Value = 625005850
$ws.Range("M89").Value = -556247.3200000001
$ws.Range("N89").Value = -625017082
# Row 98
$ws.Range("H98").Value = 47230.043
$ws.Range("I98").Value = 56521.05
$ws.Range("J98").Value = 775
$ws.Range("K98").Value = 56521.05
$ws.Range("L98").Value = 775
$ws.Range("M98").Value = -55023.05
$ws.Range("N98").Value = -3771
# Row 106
$ws.Range("H106").Value = 15657769
$ws.Range("I106").Value = 35852
$ws.Range("J106").Value = 166669630
$ws.Range("K106").Value = 35852
$ws.Range("L106").Value = 166669630
$ws.Range("M106").Value = -35221
$ws.Range("N106").Value = -166670892
# Row 107
$ws.Range("H107").Value = 672.7727
$ws.Range("I107").Value = 174.09091
$ws.Range("J107").Value = 1171.4546
$ws.Range("K107").Value = 174.09091
$ws.Range("L107").Value = 1171.4546
$ws.Range("M107").Value = 1745.90909
$ws.Range("N107").Value = -5011.4546
# Row 122
$ws.Range("H122").Value = 47230.043
$ws.Range("I122").Value = 56521.05
$ws.Range("J122").Value = 775
$ws.Range("K122").Value = 169563.15
$ws.Range("L122").Value = 2325
$ws.Range("M122").Value = -167113.15
$ws.Range("N122").Value = -7225
# Row 132
$ws.Range("H132").Value = 2343799.8
$ws.Range("I132").Value = 2978127.5
$ws.Range("K132").Value = 8934382.5
$ws.Range("M132").Value = -8931852.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 35715124
$ws.Range("I2").Value = 55556096
$ws.Range("J2").Value = 1371.5
$ws.Range("K2").Value = 55556096
$ws.Range("L2").Value = 1371.5
$ws.Range("M2").Value = -55555983
$ws.Range("N2").Value = -1597.5
# Row 26
$ws.Range("H26").Value = 3133
$ws.Range("I26").Value = 3133
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3133
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = -2803
# Row 32
$ws.Range("H32").Value = 2147.58
$ws.Range("I32").Value = 2104.6736
$ws.Range("K32").Value = 2104.6736
$ws.Range("M32").Value = -1817.6736
# Row 45
$ws.Range("H45").Value = 15152502
$ws.Range("I45").Value = 41667300
$ws.Range("J45").Value = 1188.9286
$ws.Range("K45").Value = 41667300
$ws.Range("L45").Value = 1188.9286
$ws.Range("M45").Value = -41666923
$ws.Range("N45").Value = -1942.9286
# Row 116
$ws.Range("H116").Value = 35715124
$ws.Range("I116").Value = 55556096
$ws.Range("J116").Value = 1371.5
$ws.Range("K116").Value = 55556096
$ws.Range("L116").Value = 1371.5
$ws.Range("M116").Value = -55553802
$ws.Range("N116").Value = -5959.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 35715124
$ws.Range("I3").Value = 55556096
$ws.Range("J3").Value = 1371.5
$ws.Range("K3").Value = 55556096
$ws.Range("L3").Value = 1371.5
$ws.Range("M3").Value = -55555982
$ws.Range("N3").Value = -1599.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4277460
$ws.Range("I31").Value = 3347.5652
$ws.Range("J31").Value = 10421497
$ws.Range("K31").Value = 3347.5652
$ws.Range("L31").Value = 10421497
$ws.Range("M31").Value = -3052.5652
$ws.Range("N31").Value = -10422087
# Row 34
$ws.Range("H34").Value = 4277460
$ws.Range("I34").Value = 3347.5652
$ws.Range("J34").Value = 10421497
$ws.Range("K34").Value = 3347.5652
$ws.Range("L34").Value = 10421497
$ws.Range("M34").Value = -3145.5652
$ws.Range("N34").Value = -10421901
# Row 58
$ws.Range("H58").Value = 1092.5454
$ws.Range("I58").Value = 1101.6875
$ws.Range("J58").Value = 800
$ws.Range("K58").Value = 1101.6875
$ws.Range("L58").Value = 800
$ws.Range("M58").Value = -898.6875
$ws.Range("N58").Value = -1206
# Row 136
$ws.Range("H136").Value = 1092.5454
$ws.Range("I136").Value = 1101.6875
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 3305.0625
$ws.Range("L136").Value = 2400
$ws.Range("M136").Value = -755.0625
$ws.Range("N136").Value = -7500

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 760.74
$ws.Range("I131").Value = 310
$ws.Range("J131").Value = 828.092
$ws.Range("K131").Value = 930
$ws.Range("L131").Value = 2484.276
$ws.Range("M131").Value = 4110
$ws.Range("N131").Value = -12564.276

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 47191692
$ws.Range("I70").Value = 55929984
$ws.Range("K70").Value = 55929984
$ws.Range("M70").Value = -55929714
# Row 73
$ws.Range("H73").Value = 47191692
$ws.Range("I73").Value = 55929984
$ws.Range("K73").Value = 55929984
$ws.Range("M73").Value = -55929048
# Row 102
$ws.Range("H102").Value = 1744.8
$ws.Range("I102").Value = 1235.8
$ws.Range("J102").Value = 2762.8
$ws.Range("K102").Value = 1235.8
$ws.Range("L102").Value = 2762.8
$ws.Range("M102").Value = 386.2
$ws.Range("N102").Value = -6006.8
# Row 126
$ws.Range("H126").Value = 1702.3529
$ws.Range("I126").Value = 1718.4615
$ws.Range("K126").Value = 5155.3845
$ws.Range("M126").Value = -2685.3845

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 7904.8823
$ws.Range("I132").Value = 10672.091
$ws.Range("J132").Value = 2831.6667
$ws.Range("K132").Value = 32016.273
$ws.Range("L132").Value = 8495.000100000001
$ws.Range("M132").Value = -29486.273
$ws.Range("N132").Value = -13555.0001
# Row 136
$ws.Range("H136").Value = 5033.697
$ws.Range("I136").Value = 6113.5654
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 18340.6962
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -15790.6962
$ws.Range("N136").Value = -12750

